$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    'D2' = '58.135.92'
    'E2' = '  +1.72%  '
    'D3' = '3.079.11'
    'E3' = '  +0.48%  '
    'E4' = '  +0.07%  '
    'D5' = '516.57'
    'E5' = '  +0.41%  '
    'D6' = '141.67'
    'E7' = '  +0.06%  '
    'D8' = '0.436'
    'E8' = '  +0.54%  '
    'D9' = '7.34'
    'E9' = '  +2.21%  '
    'E10' = '  -0.29%  '
    'E11' = '  +1.98%  '
    'D12' = '3.609.29'
    'E12' = '  +0.83%  '
    'D14' = '26.77'
    'E14' = '  +5.97%  '
    'D15' = '0.0000164'
    'E15' = '  -0.07%  '
    'D16' = '58.182.18'
    'D17' = '6.17'
    'E17' = '  +3.98%  '
    'D18' = '3.078.50'
    'E18' = '  +0.75%  '
    'D19' = '12.90'
    'E19' = '  -1.20%  '
    'D20' = '8.13'
    'E20' = '  -0.20%  '
    'D21' = '334.37'
    'E21' = '  -0.33%  '
    'E22' = '  -0.42%  '
    'D23' = '0.503'
    'E23' = '  +0.03%  '
    'D24' = '65.27'
    'E24' = '  +0.16%  '
    'E25' = '  +2.65%  '
    'D26' = '0.999'
    'E26' = '  -0.19%  '
    'D27' = '0.0₃0908'
    'E27' = '  -4.73%  '
    'D28' = '6.48'
    'E28' = '  +0.57%  '
    'D29' = '7.22'
    'E29' = '  +3.32%  '
    'D30' = '1.81'
    'E30' = '  +0.16%  '
    'E31' = '  +2.33%  '
    'D32' = '20.80'
    'E32' = '  +0.29%  '
    'D33' = '155.16'
    'E33' = '  +0.45%  '
    'D34' = '4.57'
    'E34' = '  +0.95%  '
    'D35' = '27.47'
    'E35' = '  +4.07%  '
    'D36' = '6.02'
    'E36' = '  +2.91%  '
    'D37' = '1.28'
    'E37' = '  +3.84%  '
    'D38' = '0.0677'
    'E38' = '  +0.99%  '
    'D39' = '3.120.40'
    'E39' = '  +0.80%  '
    'D40' = '3.91'
    'E40' = '  +2.30%  '
    'D41' = '36.60'
    'E41' = '  -0.90%  '
    'E42' = '  +0.15%  '
    'D43' = '0.658'
    'E43' = '  -1.59%  '
    'D44' = '2.277.73'
    'E44' = '  +2.04%  '
    'D45' = '0.0255'
    'E45' = '  +1.32%  '
    'E46' = '  +1.18%  '
    'D47' = '20.55'
    'E47' = '  +3.38%  '
    'D48' = '0.947'
    'E48' = '  +0.93%  '
    'D49' = '5.96'
    'E49' = '  +2.09%  '
    'D50' = '0.738'
    'E50' = '  +7.41%  '
    'D51' = '257.97'
    'E51' = '  +10.43%  '
}

foreach ($cell in $changes.Keys) {
    $range = $ws.Range($cell)
    if ($cell.Substring(0,1) -eq "D") {
        $range.NumberFormat = "@"
    }
    $range.Value = $changes[$cell]
}
